$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("writing")
$wsDash = $wb.Worksheets.Item("dashboard")

# Expand the structured table by one row, then fill in the new row's data.
$lo = $wsData.ListObjects.Item("Table1")
$lo.ListRows.Add() | Out-Null

# Copy the formatting of the prior date cell so A29 picks up the same date
# number format style as the rest of the Date column.
$wsData.Range("A28").Copy()
$wsData.Range("A29").PasteSpecial(-4122)
$wsData.Range("A29").Value = 44162

$wsData.Range("B29").Value = 236
$wsData.Range("C29").Value = 87
$wsData.Range("D29").Value = 516
$wsData.Range("E29").Value = 7562
$wsData.Range("F29").Value = 10081
$wsData.Range("G29").Value = 120
$wsData.Range("H29").Value = 117
$wsData.Range("I29").Value = 186
$wsData.Range("J29").Formula = "=SUM(B29:I29)"
$wsData.Range("K29").Formula = "=J29-J28"

# Update the dashboard chart's series so they reference the grown range.
$chartObj = $wsDash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$sDaily = $chart.SeriesCollection().Item(1)
$sDaily.Formula = "=SERIES(writing!`$K`$1,writing!`$A`$2:`$A`$29,writing!`$K`$2:`$K`$29,1)"
$sTotal = $chart.SeriesCollection().Item(2)
$sTotal.Formula = "=SERIES(writing!`$J`$1,writing!`$A`$2:`$A`$29,writing!`$J`$2:`$J`$29,2)"

# The "writing" tab becomes the active sheet/tab, with the new row selected,
# mirroring the author switching focus there after logging the day's entry.
$wsData.Range("F29").Select()
$wsData.Activate()
